# Submission Object Class remove parameter - UPDATE
# Updates the Presentation sheet's presentationDueDate, presentationStatus
# and presentationResult sample data, and refreshes the sheet
# selections / active sheet to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Presentation sheet: bulk data update
# ---------------------------------------------------------------------
$presentation = $wb.Worksheets.Item("Presentation")

# Column E (presentationDueDate) for rows 2-19 moves from 2024-04-09 to
# 2024-08-01.
for ($row = 2; $row -le 19; $row++) {
    $presentation.Cells.Item($row, 5).Value = "2024-08-01 00:00:00"
}

# Column G (presentationStatus) changes on a handful of rows away from
# the default "PENDING_BOOKING" value.
$statusOverrides = @{
    5  = "PENDING_CONFIRM"
    8  = "REJECTED"
    11 = "BOOKED"
    14 = "MARKED"
    17 = "OVERDUE"
    25 = "MARKED"
    26 = "MARKED"
    27 = "MARKED"
    28 = "MARKED"
    29 = "MARKED"
}
foreach ($row in $statusOverrides.Keys) {
    $presentation.Cells.Item($row, 7).Value = $statusOverrides[$row]
}

# Column H (presentationResult) changes from a constant "0" to an
# incrementing textual id sequence: 99, 80, 81, ..., 98, 99, 100, ..., 106
$resultValues = @("99","80","81","82","83","84","85","86","87","88","89","90","91","92","93","94","95","96","97","98","99","100","101","102","103","104","105","106")
for ($i = 0; $i -lt $resultValues.Length; $i++) {
    $row = 2 + $i
    $presentation.Cells.Item($row, 8).Value = $resultValues[$i]
}

# ---------------------------------------------------------------------
# Sheet selections / active sheet
# ---------------------------------------------------------------------
$intake = $wb.Worksheets.Item("Intake")
$intake.Activate()
$intake.Range("H27").Select()

$projectModule = $wb.Worksheets.Item("ProjectModule")
$projectModule.Activate()
$projectModule.Range("J16").Select()

$presentation.Activate()
$presentation.Range("K24").Select()

$submission = $wb.Worksheets.Item("Submission")
$submission.Activate()
$submission.Range("I12").Select()
